# Bill of materials update:
# - Raspberry Pi 4 B was actually paid for (was "already had" one, now specified as a Pi 3)
# - Added "already had" comments for the switches and dupont wire (no longer free-form blank)
# - Added an "ACTUAL COST" row total that nets out the already-owned items

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Raspberry Pi now has a real cost, comment text updated ---
$ws.Range("C2").Value = 75
$ws.Range("D2").Value = "Already had a Raspberry Pi 3"

# --- Row 5: 12*12mm push switches -- note they were already owned ---
$ws.Range("D5").Value = "Already had switches"

# --- Row 10: Dupont wire -- note it was already owned ---
$ws.Range("D10").Value = "Already had dupont wires"

# --- Row 13: totals. B13 SUM(C2:C12) stays; add an "ACTUAL COST" label + formula ---
# The old D13:E13 merge is no longer appropriate since D13/E13 now hold distinct content;
# unmerge first so both cells can carry their own value independently.
$ws.Range("D13:E13").UnMerge()
$ws.Range("D13").Value = "ACTUAL COST"
$ws.Range("E13").Formula = "=SUM(C3,C4,C6,C7,C9,C8,C11,C12)"

# Selection moved to E13 (matches where the author finished editing)
$ws.Range("E13").Select()

Write-Host "edit applied"
